# Append 15 new alumni rows (rows 10-24) to Sheet1, matching the
# "Updated member page and .xlsx sheet" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -- shared constant values reused by every new row --------------------
$img      = "../alumni/binod.jpg"
$url      = "https://istenith.com/prody/"
$position = "Development Head, Byjus"

# New alumni names, in sheet order (rows 10 .. 24)
$names = @(
    "Vivek Gusain",
    "Marmik Sharma",
    "Mayank Singh",
    "Mukul C. Mahadik",
    "Sumit Sharma",
    "Abhiraj Singh Rathore",
    "Amol Bobade",
    "Divyanshu Bhaik",
    "Kunal Kishore",
    "Sahaj Kulshrestha",
    "Parthivi Jain",
    "Varan Singh Rohila",
    "Achyut Sharma",
    "Priyanka Kumar",
    "Rishi Kumar"
)

# Per-row heights for rows 10 .. 24 (points), taken from the target sheet
$heights = @(13.8, 13.8, 13.8, 23.85, 13.8, 23.85, 13.8, 13.8, 13.8, 23.85, 13.8, 23.85, 13.8, 13.8, 13.8)

# Row 9 (the last pre-existing row) gets a slightly shorter height too
$ws.Rows.Item(9).RowHeight = 14.15

for ($i = 0; $i -lt $names.Count; $i++) {
    $r = 10 + $i

    $ws.Cells.Item($r, 1).Value = $names[$i]
    $ws.Cells.Item($r, 2).Value = 2017
    $ws.Cells.Item($r, 3).Value = 2021
    $ws.Cells.Item($r, 4).Value = $img
    $ws.Cells.Item($r, 5).Value = $url
    $ws.Cells.Item($r, 6).Value = $url
    $ws.Cells.Item($r, 7).Value = $url
    $ws.Cells.Item($r, 8).Value = $position

    # Name column wraps (new cell style: same font as default, wrapText=true)
    $ws.Cells.Item($r, 1).WrapText = $true

    $ws.Rows.Item($r).RowHeight = $heights[$i]
}

# Only the first four new rows (10-13) got live hyperlinks on E/F/G in the
# source edit -- replicate that exactly (rows 14-24 keep plain text urls).
for ($r = 10; $r -le 13; $r++) {
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 5), $url, "", "", $url) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 6), $url, "", "", $url) | Out-Null
    $ws.Hyperlinks.Add($ws.Cells.Item($r, 7), $url, "", "", $url) | Out-Null
}

# Restore the view/selection state recorded in the target sheet.
$excel.ActiveWindow.ScrollColumn = 6
$ws.Range("I22").Select()

"done"
